$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a text value into a cell while (a) preventing Excel's
# automatic date/number parsing from mangling strings that look like
# dates (e.g. "79年05月04日", "2011-11-21"), and (b) not permanently
# leaving a stray "@" number format behind - the caller re-pastes the
# correct formatting from a template cell/row right afterwards.
# ---------------------------------------------------------------------
function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# =======================================================================
# Sheet "土地" (land) - insert a new row 2 for a newly-discovered parcel,
# pushing the two existing rows down to rows 3-4.
# =======================================================================
$ws1 = $wb.Worksheets.Item("土地")

$ws1.Rows.Item(2).Insert()
$ws1.Range("A3:Q3").Copy()
$ws1.Range("A2:Q2").PasteSpecial(-4122)

$ws1.Cells.Item(2,1).Value = 15
Set-TextCell $ws1.Cells.Item(2,2) "南投縣南投市牛運堀段02670002地號"
$ws1.Cells.Item(2,3).Value = 325
Set-TextCell $ws1.Cells.Item(2,4) "全部"
Set-TextCell $ws1.Cells.Item(2,5) "廖述嘉"
Set-TextCell $ws1.Cells.Item(2,6) "79年05月04日"
Set-TextCell $ws1.Cells.Item(2,7) "共有物分割"
Set-TextCell $ws1.Cells.Item(2,8) "(超過五年）"
Set-TextCell $ws1.Cells.Item(2,9) "land"
Set-TextCell $ws1.Cells.Item(2,10) "normal"
Set-TextCell $ws1.Cells.Item(2,11) "2011-11-21"
Set-TextCell $ws1.Cells.Item(2,12) "盧秀燕"
$ws1.Cells.Item(2,13).Value = 869
Set-TextCell $ws1.Cells.Item(2,14) "tmp9eb41"
$ws1.Cells.Item(2,15).Value = 15
$ws1.Cells.Item(2,16).Value = 1
$ws1.Cells.Item(2,17).Value = 325

$ws1.Range("A3:Q3").Copy()
$ws1.Range("A2:Q2").PasteSpecial(-4122)

# =======================================================================
# Sheet "汽車" (car) - row 1 held the first (unindexed) record; add a
# proper indexed row 2 duplicating it, pushing the old row 2 to row 3.
# =======================================================================
$ws2 = $wb.Worksheets.Item("汽車")

$ws2.Rows.Item(2).Insert()
$ws2.Range("A3:G3").Copy()
$ws2.Range("A2:G2").PasteSpecial(-4122)

$ws2.Cells.Item(2,1).Value = 32
Set-TextCell $ws2.Cells.Item(2,2) "HYUNDAI"
$ws2.Cells.Item(2,3).Value = 2497
Set-TextCell $ws2.Cells.Item(2,4) "盧秀燕"
Set-TextCell $ws2.Cells.Item(2,5) "99年02月06日"
Set-TextCell $ws2.Cells.Item(2,6) "(購二手車）"
$ws2.Cells.Item(2,7).Value = 100000

$ws2.Range("A3:G3").Copy()
$ws2.Range("A2:G2").PasteSpecial(-4122)

# =======================================================================
# Sheet "存款" (deposits) - same pattern: row 1 is the unindexed first
# record; insert a proper indexed row 2 duplicating it.
# =======================================================================
$ws3 = $wb.Worksheets.Item("存款")

$ws3.Rows.Item(2).Insert()
$ws3.Range("A3:F3").Copy()
$ws3.Range("A2:F2").PasteSpecial(-4122)

$ws3.Cells.Item(2,1).Value = 47
Set-TextCell $ws3.Cells.Item(2,2) "臺灣銀行群賢分行"
Set-TextCell $ws3.Cells.Item(2,3) "活期儲蓄存款"
Set-TextCell $ws3.Cells.Item(2,4) "新臺幣"
Set-TextCell $ws3.Cells.Item(2,5) "盧秀燕"
$ws3.Cells.Item(2,6).Value = 4752062

$ws3.Range("A3:F3").Copy()
$ws3.Range("A2:F2").PasteSpecial(-4122)

# =======================================================================
# Sheet "保險" (insurance) - same pattern again.
# =======================================================================
$ws4 = $wb.Worksheets.Item("保險")

$ws4.Rows.Item(2).Insert()
$ws4.Range("A3:E3").Copy()
$ws4.Range("A2:E2").PasteSpecial(-4122)

$ws4.Cells.Item(2,1).Value = 87
Set-TextCell $ws4.Cells.Item(2,2) "南山人壽"
Set-TextCell $ws4.Cells.Item(2,3) "子女教育保險"
Set-TextCell $ws4.Cells.Item(2,4) "廖述嘉"
Set-TextCell $ws4.Cells.Item(2,5) "保險期間：951811718(22年)年繳保費應繳`$71400"

$ws4.Range("A3:E3").Copy()
$ws4.Range("A2:E2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
